$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 66666868
$ws.Range("I5").Value = 90909180
$ws.Range("J5").Value = 499.25
$ws.Range("K5").Value = 90909180
$ws.Range("L5").Value = 499.25
$ws.Range("M5").Value = -90909065
$ws.Range("N5").Value = -729.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 973.56177
$ws.Range("I129").Value = 327.14285
$ws.Range("J129").Value = 1028.7439
$ws.Range("K129").Value = 981.4285500000001
$ws.Range("L129").Value = 3086.2317
$ws.Range("M129").Value = 4018.57145
$ws.Range("N129").Value = -13086.2317

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5230716
$ws.Range("I137").Value = 7354105
$ws.Range("J137").Value = 73913.5
$ws.Range("K137").Value = 22062315
$ws.Range("L137").Value = 221740.5
$ws.Range("M137").Value = -22059765
$ws.Range("N137").Value = -226840.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6825
$ws.Range("I61").Value = 1800
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 1800
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -1588
$ws.Range("N61").Value = -8924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4676.9644
$ws.Range("I74").Value = 5038.44
$ws.Range("J74").Value = 1664.6666
$ws.Range("K74").Value = 5038.44
$ws.Range("L74").Value = 1664.6666
$ws.Range("M74").Value = -4164.44
$ws.Range("N74").Value = -3412.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4676.9644
$ws.Range("I77").Value = 5038.44
$ws.Range("J77").Value = 1664.6666
$ws.Range("K77").Value = 25192.2
$ws.Range("L77").Value = 8323.333000000001
$ws.Range("M77").Value = -20824.2
$ws.Range("N77").Value = -17059.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1037.25
$ws.Range("I97").Value = 1050
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = -554
$ws.Range("N97").Value = -1991

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6825
$ws.Range("I136").Value = 1800
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 5400
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -2850
$ws.Range("N136").Value = -30600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.2
$ws.Range("I7").Value = 59
$ws.Range("J7").Value = 87
$ws.Range("K7").Value = 59
$ws.Range("L7").Value = 87
$ws.Range("M7").Value = 54
$ws.Range("N7").Value = -313

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2624.12
$ws.Range("I58").Value = 1727.6154
$ws.Range("J58").Value = 3595.3333
$ws.Range("K58").Value = 1727.6154
$ws.Range("L58").Value = 3595.3333
$ws.Range("M58").Value = -1524.6154
$ws.Range("N58").Value = -4001.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2624.12
$ws.Range("I136").Value = 1727.6154
$ws.Range("J136").Value = 3595.3333
$ws.Range("K136").Value = 5182.8462
$ws.Range("L136").Value = 10785.9999
$ws.Range("M136").Value = -2632.8462
$ws.Range("N136").Value = -15885.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 34.5
$ws.Range("I14").Value = 34.5
$ws.Range("K14").Value = 103.5
$ws.Range("M14").Value = 69.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1437354.6
$ws.Range("I113").Value = 3135181.8
$ws.Range("J113").Value = 731.53845
$ws.Range("K113").Value = 9405545.399999999
$ws.Range("L113").Value = 2194.61535
$ws.Range("M113").Value = -9403375.399999999
$ws.Range("N113").Value = -6534.61535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 3886.625
$ws.Range("I125").Value = 1030
$ws.Range("J125").Value = 4838.8335
$ws.Range("K125").Value = 3090
$ws.Range("L125").Value = 14516.5005
$ws.Range("M125").Value = 1830
$ws.Range("N125").Value = -24356.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 791.3086499999999
$ws.Range("J131").Value = 901.5484
$ws.Range("L131").Value = 2704.6452
$ws.Range("N131").Value = -12784.6452

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1684818.4
$ws.Range("I132").Value = 1350
$ws.Range("J132").Value = 2526552.5
$ws.Range("K132").Value = 12150
$ws.Range("L132").Value = 22738972.5
$ws.Range("M132").Value = -9620
$ws.Range("N132").Value = -22744032.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6668466
$ws.Range("I122").Value = 9092436
$ws.Range("J122").Value = 2547
$ws.Range("K122").Value = 27277308
$ws.Range("L122").Value = 7641
$ws.Range("M122").Value = -27274858
$ws.Range("N122").Value = -12541

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 36000
$ws.Range("I61").Value = 51500
$ws.Range("K61").Value = 51500
$ws.Range("M61").Value = -51298

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3211.4285
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 3696
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 3696
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -5194

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3211.4285
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 3696
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 18480
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -25968

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1585.5714
$ws.Range("I82").Value = 1585.5714
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1585.5714
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1224.5714
$ws.Range("N82").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1585.5714
$ws.Range("I85").Value = 1585.5714
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1585.5714
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -337.5714
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 36000
$ws.Range("I113").Value = 51500
$ws.Range("K113").Value = 51500
$ws.Range("M113").Value = -49330

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2946.8928
$ws.Range("I122").Value = 2683.9473
$ws.Range("J122").Value = 3502
$ws.Range("K122").Value = 8051.841899999999
$ws.Range("L122").Value = 10506
$ws.Range("M122").Value = -5601.841899999999
$ws.Range("N122").Value = -15406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1711.1464
$ws.Range("I136").Value = 1402.0938
$ws.Range("K136").Value = 4206.2814
$ws.Range("M136").Value = -1656.2814
